# Publish IG 1.0.1
# - Identifier: drop the "id: " prefix
# - Version: 1.0.0 -> 1.0.1
# - Contact: replace placeholder with the real MedCom contact
# - Insert a new "Jurisdiction" metadata row (blank value) right after "Contact",
#   pushing Description..Count down by one row
# - Supplements count is still current, sheet2 (Concepts) untouched

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Make room for the new "Jurisdiction" row: shift rows 12..22 down to 13..23 ---
# Walk bottom-up so we never overwrite a row before it has been copied down.
for ($r = 22; $r -ge 12; $r--) {
    $destRow = $r + 1
    $aVal = $ws.Cells.Item($r, 1).Text
    $bVal = $ws.Cells.Item($r, 2).Text

    # Leading apostrophe forces text so Excel doesn't auto-coerce things like
    # "false"/"true"/"6"/dates into bool/number/date cells.
    $ws.Cells.Item($destRow, 1).Value = "'" + $aVal
    $ws.Cells.Item($destRow, 2).Value = "'" + $bVal
}

# --- Write the new Jurisdiction row (blank value) into the freed-up row 12 ---
$ws.Cells.Item(12, 1).Value = "'Jurisdiction"
$ws.Cells.Item(12, 2).Value = "'"

# --- Field value updates ---
$ws.Cells.Item(3, 2).Value = "'1.2.208.184.100.1"
$ws.Cells.Item(4, 2).Value = "'1.0.1"
$ws.Cells.Item(11, 2).Value = "'MedCom (http://www.medcom.dk)"

Write-Output "done"
